$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.586.57"
$ws.Range("E2").Value = "  -1.02%  "
$ws.Range("D3").Value = "1.597.98"
$ws.Range("E3").Value = "  -1.67%  "
$ws.Range("E4").Value = "  +0.24%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "208.92"
$ws.Range("E5").Value = "  -0.99%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.500"
$ws.Range("E6").Value = "  -4.18%  "
$ws.Range("E7").Value = "  +0.32%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "22.42"
$ws.Range("E8").Value = "  -4.25%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.252"
$ws.Range("E9").Value = "  -1.89%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0590"
$ws.Range("E10").Value = "  -3.47%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0873"
$ws.Range("E11").Value = "  -0.63%  "
$ws.Range("D12").Value = "1.825.52"
$ws.Range("E12").Value = "  -1.60%  "
$ws.Range("D13").Value = "1.632.48"
$ws.Range("E13").Value = "  +1.01%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.86"
$ws.Range("E14").Value = "  -3.91%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.542"
$ws.Range("E15").Value = "  -3.39%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "63.66"
$ws.Range("E16").Value = "  -2.53%  "
$ws.Range("D17").Value = "27.574.20"
$ws.Range("E17").Value = "  -1.02%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "218.03"
$ws.Range("E18").Value = "  -4.96%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.44"
$ws.Range("E19").Value = "  -2.69%  "
$ws.Range("D20").Value = "0.0₃0691"
$ws.Range("E20").Value = "  -4.21%  "
$ws.Range("E21").Value = "  +0.26%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.20"
$ws.Range("E22").Value = "  -2.66%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.76"
$ws.Range("E23").Value = "  -3.74%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.02"
$ws.Range("E24").Value = "  -0.34%  "
$ws.Range("E25").Value = "  -0.43%  "
$ws.Range("B26").Value = "BinanceUSD"
$ws.Range("C26").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("E26").Value = "  +0.26%  "
$ws.Range("B27").Value = "Cosmos"
$ws.Range("C27").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "6.72"
$ws.Range("E27").Value = "  -2.61%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.03"
$ws.Range("E28").Value = "  -3.02%  "
$ws.Range("E29").Value = "  -4.53%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.16"
$ws.Range("E30").Value = "  -1.54%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0467"
$ws.Range("E31").Value = "  -2.78%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.31"
$ws.Range("E32").Value = "  -2.90%  "
$ws.Range("D33").Value = "1.370.94"
$ws.Range("E33").Value = "  -1.55%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.96"
$ws.Range("E34").Value = "  -3.52%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.55"
$ws.Range("E35").Value = "  -1.03%  "
$ws.Range("E36").Value = "  -5.14%  "
$ws.Range("E37").Value = "  -0.79%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0165"
$ws.Range("E38").Value = "  -2.66%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.536"
$ws.Range("E39").Value = "  -3.37%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.816"
$ws.Range("E40").Value = "  -4.86%  "
$ws.Range("E41").Value = "  +0.32%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.979"
$ws.Range("E42").Value = "  -4.51%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "63.89"
$ws.Range("E43").Value = "  -2.67%  "
$ws.Range("E44").Value = "  -3.49%  "
$ws.Range("E45").Value = "  -3.67%  "
$ws.Range("D46").Value = "1.735.43"
$ws.Range("E46").Value = "  -1.72%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.08"
$ws.Range("E47").Value = "  -3.85%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "87.97"
$ws.Range("B49").Value = "BabyDogeCoin"
$ws.Range("C49").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D49").Value = "0.0₇0999"
$ws.Range("E49").Value = "  -5.02%  "
$ws.Range("B50").Value = "Algorand"
$ws.Range("C50").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0970"
$ws.Range("E50").Value = "  -4.51%  "
$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0497"
$ws.Range("E51").Value = "  -1.19%  "
